$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# like "64.058.15" or "0.540" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.058.15"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").Value = "3.152.82"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "591.90"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("D6").Value = "147.36"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.144.90"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +22.05%  "
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  +5.66%  "
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  +10.44%  "
$ws.Range("D14").Value = "36.02"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "3.675.86"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.984.92"
$ws.Range("E17").Value = "  +4.72%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.153.55"
$ws.Range("E18").Value = "  +3.72%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "7.18"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "471.56"
$ws.Range("E20").Value = "  +5.73%  "
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +4.54%  "
$ws.Range("D24").Value = "13.42"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "82.85"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "8.67"
$ws.Range("E27").Value = "  +7.24%  "
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").Value = "6.88"
$ws.Range("E31").Value = "  +6.90%  "
$ws.Range("D32").Value = "27.22"
$ws.Range("E32").Value = "  +3.04%  "
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("D34").Value = "0.0₃0881"
$ws.Range("E34").Value = "  +8.84%  "
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +12.28%  "
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("D37").Value = "3.43"
$ws.Range("E37").Value = "  +16.58%  "
$ws.Range("D38").Value = "6.17"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "50.88"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").Value = "448.90"
$ws.Range("E40").Value = "  +10.18%  "
$ws.Range("D41").Value = "8.73"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "0.0378"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").Value = "2.923.15"
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("D44").Value = "0.281"
$ws.Range("E44").Value = "  +7.60%  "
$ws.Range("D45").Value = "0.113"
$ws.Range("E45").Value = "  +5.12%  "
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").Value = "125.80"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "24.94"
$ws.Range("E50").Value = "  +3.99%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "34.28"
$ws.Range("E51").Value = "  -8.08%  "

# Restore the default cell style on column D so no stray number-format
# style index is left referenced on these cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
